$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = [char]0x2022

$lines = @(
  "$bullet The amendment introduces changes to the definition of ""controlling interest"" for Real Estate Investment Trusts (REITs), depending on whether their shares are listed or not.",
  "$bullet If listed, the definition will align with regulations under Section 11 of the SEBI Act and the Companies Act, 2013. If unlisted, it will follow the Companies Act, 2013 definition.",
  "$bullet The amendment expands the definition of ""common infrastructure"" to include facilities such as power plants, water treatment plants, waste treatment plants, and other amenities that exclusively supply or cater to REIT, its HoldCo(s) or SPV(s).",
  "$bullet Any excess production or capacity from common infrastructure can now be sold or supplied to a central or state grid or utility, subject to specific conditions.",
  "$bullet The amendment adds a new definition for ""employee unit option scheme,"" which is a scheme under which the manager grants unit options to its employees through an employee benefit trust.",
  "$bullet Employees of the manager will include all directors, except independent directors."
)

$newText = [string]::Join([char]10, $lines)

$ws.Range("J2").Value = $newText
